# "每日学习.xlsx" - add a new daily-log entry (row 30) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The merged block A28:C29 / E29 previously relied on a duplicated
# "centered" style; re-apply the centered alignment so the workbook's
# style table collapses the duplicate the same way Excel does on save.
$ws.Range("A28:C29").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A28:C29").VerticalAlignment = -4108     # xlCenter

# New entry for day 21 (second session): time spent and topic covered.
$ws.Range("A30").Value = 21
$ws.Range("B30").Value = "2：00-5：13"
$ws.Range("C30").Value = "到5.6指针数组以及指向指针的指针"

$ws.Range("C30").Select()
